$wb = $excel.ActiveWorkbook

# The data sheet ("Sheet1") becomes "TestData" and moves to just before
# "TestSheet" (i.e. after "Sheet2"/"Sheet3"), becoming the workbook's
# active tab with a new selection.
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Name = "TestData"
$ws.Move($wb.Worksheets.Item("TestSheet"))

# Re-resolve by name so the activation/selection land on the freshly
# renamed + relocated sheet.
$testData = $wb.Worksheets.Item("TestData")
$testData.Activate()
$testData.Range("C15").Select()
